$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.741.89'
$ws.Range('E2').Value = '  -0.32%  '
$ws.Range('D3').Value = '1.634.33'
$ws.Range('E3').Value = '  -0.42%  '
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '215.35'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.24%  '
$ws.Range('E6').Value = '  -0.67%  '
$ws.Range('E7').Value = '  -0.14%  '
$ws.Range('E8').Value = '  -0.26%  '
$ws.Range('E9').Value = '  -1.21%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.57'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -4.23%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0786'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.71%  '
$ws.Range('E12').Value = '  -0.65%  '
$ws.Range('D13').Value = '1.859.84'
$ws.Range('E13').Value = '  -0.47%  '
$ws.Range('D14').Value = '1.633.57'
$ws.Range('E14').Value = '  -0.64%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.555'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.29%  '
$ws.Range('E16').Value = '  -0.13%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '62.75'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.08%  '
$ws.Range('D18').Value = '25.769.96'
$ws.Range('E19').Value = '  -0.19%  '
$ws.Range('E20').Value = '  +1.07%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '193.84'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.48%  '
$ws.Range('E22').Value = '  -0.04%  '
$ws.Range('E24').Value = '  -0.13%  '
$ws.Range('E25').Value = '  +2.51%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '140.44'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.72%  '
$ws.Range('E27').Value = '  -1.49%  '
$ws.Range('E28').Value = '  +0.66%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.52'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.69%  '
$ws.Range('E30').Value = '  -0.19%  '
$ws.Range('E31').Value = '  -0.47%  '
$ws.Range('E32').Value = '  +1.19%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.24'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.33%  '
$ws.Range('E34').Value = '  +0.73%  '
$ws.Range('E35').Value = '  +0.62%  '
$ws.Range('E36').Value = '  -0.58%  '
$ws.Range('B37').Value = 'Maker'
$ws.Range('C37').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D37').Value = '1.120.79'
$ws.Range('E37').Value = '  -1.19%  '
$ws.Range('B38').Value = 'MXToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.52'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.85%  '
$ws.Range('B39').Value = 'ImmutableX'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.548'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E40').Value = '  -1.14%  '
$ws.Range('E41').Value = '  +0.51%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.58'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.88%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '99.57'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.62%  '
$ws.Range('E44').Value = '  +0.12%  '
$ws.Range('D45').Value = '1.768.99'
$ws.Range('D46').Value = '0.0₆0110'
$ws.Range('E46').Value = '  -3.56%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '55.03'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.06%  '
$ws.Range('E48').Value = '  -2.32%  '
$ws.Range('E49').Value = '  -0.53%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.57'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -3.24%  '
$ws.Range('B51').Value = 'Frax'
$ws.Range('C51').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.00'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.50%  '
